$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.061.97"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.718.67"
$ws.Range("E3").Value = "  -1.54%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.64%  "

# Row 5 - BNB
$ws.Range("D5").Value = "598.94"
$ws.Range("E5").Value = "  +3.57%  "

# Row 6 - Solana
$ws.Range("D6").Value = "184.74"
$ws.Range("E6").Value = "  +13.74%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.714.70"
$ws.Range("E7").Value = "  -6.52%  "

# Row 8 - XRP
$ws.Range("D8").Value = "0.633"
$ws.Range("E8").Value = "  -2.72%  "

# Row 9 - USDC
$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  -0.05%  "

# Row 10 - Cardano
$ws.Range("D10").Value = "0.723"
$ws.Range("E10").Value = "  -0.58%  "

# Row 11 - Dogecoin
$ws.Range("D11").Value = "0.163"
$ws.Range("E11").Value = "  -4.32%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "56.65"
$ws.Range("E12").Value = "  +10.00%  "

# Row 13 - ShibaInu
$ws.Range("D13").Value = "0.0000293"
$ws.Range("E13").Value = "  -5.81%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "10.75"
$ws.Range("E14").Value = "  -1.66%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.346.86"
$ws.Range("E15").Value = "  -0.62%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.756.08"
$ws.Range("E16").Value = "  -0.31%  "

# Row 17 - Chainlink
$ws.Range("D17").Value = "19.45"
$ws.Range("E17").Value = "  -3.90%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -1.59%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "12.91"
$ws.Range("E19").Value = "  -4.17%  "

# Row 20 - Polygon
$ws.Range("E20").Value = "  -4.53%  "

# Row 21 - WrappedBTC
$ws.Range("D21").Value = "69.183.43"
$ws.Range("E21").Value = "  +0.21%  "

# Row 22 - BitcoinCash
$ws.Range("D22").Value = "412.03"
$ws.Range("E22").Value = "  -3.10%  "

# Row 23 - PancakeSwap
$ws.Range("D23").Value = "4.64"
$ws.Range("E23").Value = "  +0.99%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "89.21"
$ws.Range("E24").Value = "  -2.25%  "

# Row 25 - ImmutableX
$ws.Range("D25").Value = "3.05"
$ws.Range("E25").Value = "  -4.74%  "

# Rows 26-28 reordered: InternetComputer(DFINITY), RenderToken, Toncoin
#                   ->  Toncoin, InternetComputer(DFINITY), RenderToken
$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "4.12"
$ws.Range("E26").Value = "  +5.55%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "12.88"
$ws.Range("E27").Value = "  -4.38%  "

$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "10.95"
$ws.Range("E28").Value = "  -0.21%  "

# Row 29 - LEO
$ws.Range("D29").Value = "6.09"
$ws.Range("E29").Value = "  +3.30%  "

# Row 30 - Filecoin
$ws.Range("D30").Value = "9.50"
$ws.Range("E30").Value = "  -6.18%  "

# Row 31 - EthereumClassic
$ws.Range("D31").Value = "32.93"
$ws.Range("E31").Value = "  -3.44%  "

# Row 32 - NEARProtocol
$ws.Range("D32").Value = "7.34"
$ws.Range("E32").Value = "  -7.55%  "

# Row 33 - Cosmos
$ws.Range("D33").Value = "12.49"
$ws.Range("E33").Value = "  -4.98%  "

# Row 34 - Hedera
$ws.Range("E34").Value = "  -2.70%  "

# Rows 35-36 reordered: InjectiveProtocol, Bittensor -> Bittensor, InjectiveProtocol
$ws.Range("B35").Value = "Bittensor"
$ws.Range("C35").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D35").Value = "614.46"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("B36").Value = "InjectiveProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D36").Value = "43.85"
$ws.Range("E36").Value = "  -6.40%  "

# Row 37 - OKB
$ws.Range("D37").Value = "65.25"
$ws.Range("E37").Value = "  -4.15%  "

# Row 38 - PEPE
$ws.Range("D38").Value = "0.0₃0900"
$ws.Range("E38").Value = "  -5.93%  "

# Row 39 - TheGraph
$ws.Range("D39").Value = "0.404"
$ws.Range("E39").Value = "  -2.62%  "

# Row 40 - Dai
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.26%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  +0.57%  "

# Row 42 - Kaspa
$ws.Range("D42").Value = "0.137"
$ws.Range("E42").Value = "  -2.89%  "

# Row 43 - ThetaToken
$ws.Range("D43").Value = "3.05"
$ws.Range("E43").Value = "  -3.05%  "

# Row 44 - Fetch.AI
$ws.Range("D44").Value = "2.73"
$ws.Range("E44").Value = "  -1.85%  "

# Row 45 - dogwifhat
$ws.Range("D45").Value = "2.98"
$ws.Range("E45").Value = "  -4.23%  "

# Row 46 - VeChain
$ws.Range("D46").Value = "0.0444"
$ws.Range("E46").Value = "  -2.67%  "

# Row 47 - THORChain
$ws.Range("D47").Value = "9.30"
$ws.Range("E47").Value = "  -3.43%  "

# Row 48 - Stellar
$ws.Range("D48").Value = "0.135"
$ws.Range("E48").Value = "  -3.33%  "

# Rows 49-51 reordered: Maker, ApeXProtocol, WEMIXToken -> WEMIXToken, Maker, ApeXProtocol
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.73"
$ws.Range("E49").Value = "  -0.89%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.784.89"
$ws.Range("E50").Value = "  -0.33%  "

$ws.Range("B51").Value = "ApeXProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D51").Value = "3.20"
$ws.Range("E51").Value = "  -0.62%  "
